# Update bracket-round winner names on the "Bracket" sheet.
# These are the round-2/round-3/etc. "winner" cells whose value is
# re-picked (animal-name text, stored as a shared string) as the
# single-elimination bracket progresses.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bracket")

$ws.Range("M4").Value = 'Chestnut-headed Bee-eater'
$ws.Range("D6").Value = 'Raven '
$ws.Range("N6").Value = 'Chestnut-headed Bee-eater'
$ws.Range("F8").Value = 'Sperm Whale '
$ws.Range("L8").Value = 'Chestnut-headed Bee-eater'
$ws.Range("D10").Value = 'Starling'
$ws.Range("N10").Value = 'Wichita Mountains Pillsnail'
$ws.Range("E12").Value = 'Starling'
$ws.Range("M12").Value = 'Koala'
$ws.Range("N14").Value = 'Koala'
$ws.Range("G16").Value = 'Sperm Whale '
$ws.Range("K16").Value = 'Great Skua'
$ws.Range("D18").Value = 'Boar'
$ws.Range("N18").Value = 'Fork-marked Lemur'
$ws.Range("E20").Value = 'Boar'
$ws.Range("M20").Value = 'Fork-marked Lemur'
$ws.Range("N22").Value = 'Velvet Worm'
$ws.Range("F24").Value = 'Bigeye Houndshark'
$ws.Range("D26").Value = 'Bigeye Houndshark'
$ws.Range("E28").Value = 'Bigeye Houndshark'
$ws.Range("D30").Value = 'Tarzan Chameleon'
$ws.Range("H32").Value = 'Sperm Whale '
$ws.Range("I32").Value = 'Sperm Whale '
$ws.Range("J32").Value = 'Great Skua'
$ws.Range("D34").Value = 'Painted Wild Dog'
$ws.Range("E36").Value = 'Painted Wild Dog'
$ws.Range("D38").Value = 'Himalayan Monal'
$ws.Range("F40").Value = 'Marbled Polecat'
$ws.Range("N42").Value = 'Wrinkle-faced Bat'
$ws.Range("M44").Value = 'Wrinkle-faced Bat'
$ws.Range("G48").Value = 'Marbled Polecat'
$ws.Range("K48").Value = 'Northern Elephant Seal'
$ws.Range("D50").Value = 'Peacock Mantis Shrimp '
$ws.Range("N50").Value = 'Howler Monkey '
$ws.Range("M52").Value = 'Howler Monkey '
$ws.Range("N54").Value = 'Elegant Dancing Frog'
$ws.Range("L56").Value = 'Howler Monkey '
$ws.Range("N58").Value = 'Flame Bowerbird'
$ws.Range("M60").Value = 'Honey Bee'
$ws.Range("N62").Value = 'Honey Bee'
